$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12: Don't Be So Tallow | Beeswax
$ws.Range("H12").Value = 125
$ws.Range("I12").Value = 125
$ws.Range("K12").Value = 125
$ws.Range("M12").Value = 45

# Row 28: The Writing Is Not on the Wall | Enchanted Silver Ink
$ws.Range("H28").Value = 2284.1428
$ws.Range("I28").Value = 2514.8333
$ws.Range("J28").Value = 900
$ws.Range("K28").Value = 2514.8333
$ws.Range("L28").Value = 900
$ws.Range("M28").Value = -2029.8333
$ws.Range("N28").Value = -1870

# Row 100: Asking for a Friend | Beetle Glue
$ws.Range("H100").Value = 4248.25
$ws.Range("I100").Value = 4248.25
$ws.Range("K100").Value = 4248.25
$ws.Range("M100").Value = -3707.25

# Row 107: Another Man's Ink | Enchanted Truegold Ink
$ws.Range("H107").Value = 660.7692
$ws.Range("I107").Value = 660.7692
$ws.Range("K107").Value = 660.7692
$ws.Range("M107").Value = 1259.2308

# Row 111: An Eye for Healing | Grade 1 Dexterity Alkahest
$ws.Range("H111").Value = 1743
$ws.Range("I111").Value = 1914.5
$ws.Range("J111").Value = 1400
$ws.Range("K111").Value = 5743.5
$ws.Range("L111").Value = 4200
$ws.Range("M111").Value = -2676.5
$ws.Range("N111").Value = -10334

# Row 116: Growing Up | Growth Formula Kappa
$ws.Range("H116").Value = 4803
$ws.Range("I116").Value = 4801
$ws.Range("K116").Value = 4801
$ws.Range("M116").Value = -1359

# Row 133: Big Brush, Big Dreams | Ginseng Angle Brush
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# Row 134: Binding Spells | Crocodileskin Index
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# Row 135: For Tired Minds | Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 2976.5557
$ws.Range("I135").Value = 2976.5557
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 26789.0013
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -24254.0013
$ws.Range("N135").ClearContents()

# Row 136: I Like Big Brush and I Cannot Lie | Dark Mahogany Round Brush
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 3509.353
$ws.Range("I137").Value = 707.375
$ws.Range("K137").Value = 2122.125
$ws.Range("M137").Value = 427.875

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 4148.316
$ws.Range("I138").Value = 2790.875
$ws.Range("J138").Value = 4510.3
$ws.Range("K138").Value = 8372.625
$ws.Range("L138").Value = 13530.9
$ws.Range("M138").Value = -3232.625
$ws.Range("N138").Value = -23810.9

$ws = $wb.Worksheets.Item("ARM")
# Row 13: Get into Their Heads | Bronze Chain Coif
$ws.Range("H13").Value = 50000750
$ws.Range("I13").Value = 50000750
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 50000750
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("M13").Value = -50000606

# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 6470.2
$ws.Range("I32").Value = 5244.722
$ws.Range("K32").Value = 5244.722
$ws.Range("M32").Value = -4957.722

# Row 45: Hollow Hallmarks | Mythril Ingot
$ws.Range("H45").Value = 2209.2354
$ws.Range("I45").Value = 1955.7
$ws.Range("J45").Value = 2571.4285
$ws.Range("K45").Value = 1955.7
$ws.Range("L45").Value = 2571.4285
$ws.Range("M45").Value = -1578.7
$ws.Range("N45").Value = -3325.4285

# Row 133: Shielding My Students | Mountain Chromite Tower Shield
$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -85060

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt | Iron Ingot
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

# Row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 1949.6666
$ws.Range("I86").Value = 1949.6666
$ws.Range("K86").Value = 1949.6666
$ws.Range("M86").Value = -826.6666

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 1949.6666
$ws.Range("I89").Value = 1949.6666
$ws.Range("K89").Value = 9748.333000000001
$ws.Range("M89").Value = -4132.333000000001

# Row 107: The Gold Experience | Deepgold Nugget
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 2762.4546
$ws.Range("I134").Value = 2762.4546
$ws.Range("K134").Value = 8287.363799999999
$ws.Range("M134").Value = -5752.363799999999

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent | Maple Lumber
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

# Row 22: Driving Up the Wall | Elm Lumber
$ws.Range("H22").Value = 2250
$ws.Range("I22").Value = 2333.3333
$ws.Range("K22").Value = 2333.3333
$ws.Range("M22").Value = -1983.3333

# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 6405.5
$ws.Range("J31").Value = 17007.2
$ws.Range("L31").Value = 17007.2
$ws.Range("N31").Value = -17597.2

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 6405.5
$ws.Range("J34").Value = 17007.2
$ws.Range("L34").Value = 17007.2
$ws.Range("N34").Value = -17411.2

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 2180
$ws.Range("I58").Value = 2180
$ws.Range("K58").Value = 2180
$ws.Range("M58").Value = -1977

# Row 107: Built to Last | White Oak Lumber
$ws.Range("H107").Value = 2780.6667
$ws.Range("J107").Value = 3849
$ws.Range("L107").Value = 3849
$ws.Range("N107").Value = -7689

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 3440.0625
$ws.Range("I132").Value = 2849.6155
$ws.Range("J132").Value = 5998.6665
$ws.Range("K132").Value = 8548.8465
$ws.Range("L132").Value = 17995.9995
$ws.Range("M132").Value = -6018.8465
$ws.Range("N132").Value = -23055.9995

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 4622.636
$ws.Range("I134").Value = 4622.636
$ws.Range("K134").Value = 13867.908
$ws.Range("M134").Value = -11332.908

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 2180
$ws.Range("I136").Value = 2180
$ws.Range("K136").Value = 6540
$ws.Range("M136").Value = -3990

$ws = $wb.Worksheets.Item("CUL")
# Row 60: Drinking to Your Health | Mulled Tea
$ws.Range("H60").Value = 726.8570999999999
$ws.Range("I60").Value = 756.8333
$ws.Range("J60").Value = 547
$ws.Range("K60").Value = 2270.4999
$ws.Range("L60").Value = 1641
$ws.Range("M60").Value = -2019.4999
$ws.Range("N60").Value = -2143

# Row 107: Slippery Service | Frantoio Oil
$ws.Range("H107").Value = 2764
$ws.Range("J107").Value = 582.5
$ws.Range("L107").Value = 1747.5
$ws.Range("N107").Value = -5587.5

# Row 113: Can't Eat Just One | Night Vinegar
$ws.Range("H113").Value = 659.6667
$ws.Range("I113").Value = 432.7143
$ws.Range("J113").Value = 977.4
$ws.Range("K113").Value = 1298.1429
$ws.Range("L113").Value = 2932.2
$ws.Range("M113").Value = 871.8571000000002
$ws.Range("N113").Value = -7272.2

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 1431.2
$ws.Range("I131").Value = 990.5
$ws.Range("J131").Value = 1499
$ws.Range("K131").Value = 2971.5
$ws.Range("L131").Value = 4497
$ws.Range("M131").Value = 2068.5
$ws.Range("N131").Value = -14577

# Row 132: More Mezcal | Cooking Mezcal
$ws.Range("H132").Value = 3846.0908
$ws.Range("I132").Value = 1969.2
$ws.Range("K132").Value = 17722.8
$ws.Range("M132").Value = -15192.8

# Row 139: Najoothie | Wild Banana Blend
$ws.Range("H139").Value = 7000
$ws.Range("J139").Value = 7000
$ws.Range("L139").Value = 21000
$ws.Range("N139").Value = -31280

$ws = $wb.Worksheets.Item("GSM")
# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 3840.3572
$ws.Range("I132").Value = 3276.9
$ws.Range("J132").Value = 5249
$ws.Range("K132").Value = 9830.700000000001
$ws.Range("L132").Value = 15747
$ws.Range("M132").Value = -7300.700000000001
$ws.Range("N132").Value = -20807

$ws = $wb.Worksheets.Item("LTW")
# Row 68: You Could Say It's a Moving Target | Wyvern Leather
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

# Row 71: They Call It Bloody Mary (L) | Wyvern Leather
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 5216.1665
$ws.Range("I132").Value = 4768
$ws.Range("J132").Value = 5664.3335
$ws.Range("K132").Value = 14304
$ws.Range("L132").Value = 16993.0005
$ws.Range("M132").Value = -11774
$ws.Range("N132").Value = -22053.0005

# Row 133: The Perfect Accessory | Loboskin Amulet of Fending
$ws.Range("H133").Value = 135000
$ws.Range("J133").Value = 135000
$ws.Range("L133").Value = 135000
$ws.Range("N133").Value = -140060

# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Range("H136").Value = 38580.082
$ws.Range("I136").Value = 15742.25
$ws.Range("K136").Value = 47226.75
$ws.Range("M136").Value = -44676.75

$ws = $wb.Worksheets.Item("WVR")
# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 1990.8572
$ws.Range("I126").Value = 1809.25
$ws.Range("J126").Value = 2233
$ws.Range("K126").Value = 5427.75
$ws.Range("L126").Value = 6699
$ws.Range("M126").Value = -2957.75
$ws.Range("N126").Value = -11639

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 2987.5
$ws.Range("I132").Value = 2645.44
$ws.Range("J132").Value = 4697.8
$ws.Range("K132").Value = 7936.32
$ws.Range("L132").Value = 14093.4
$ws.Range("M132").Value = -5406.32
$ws.Range("N132").Value = -19153.4

# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 8806.777
$ws.Range("I136").Value = 8806.777
$ws.Range("K136").Value = 26420.331
$ws.Range("M136").Value = -23870.331
